# Update the workbook's version/build string wherever it appears.
$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Range("A2").Value = "Version: $newVersion"
$aboutSheet.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for KWK Budryk Coal Mine, Poland, M1288, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")
for ($row = 2; $row -le 10; $row++) {
    $cell = $dataSheet.Range("S$row")
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
